$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the hyperlink that currently lives on C4 and restore that
#     cell to the default ("Normal") style, then drop the now unused
#     "Hyperlink" cell style from the workbook's style table. ---
$ws.Hyperlinks.Delete()
$ws.Range("C4").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

# --- Add the two new log rows. Values are entered URL-cell first so the
#     shared-string table ends up in the same order as the target file. ---
$ws.Range("C6").Value = "https://www.interviewquery.com/questions/like-tracker"
$ws.Range("A6").Value = "Like tracker"
$ws.Range("B6").Value = "Easy"
$ws.Range("D6").Value = "case to specify a date in a timestamp column - alternatives - timestamp can be cast as date usiong cast(col as date) or using specifier :::date or date format)"
$ws.Rows.Item(6).RowHeight = 68

$ws.Range("C7").Value = "https://www.interviewquery.com/questions/manager-team-sizes"
$ws.Range("A7").Value = "Manager team sizes"
$ws.Range("B7").Value = "Easy"
$ws.Range("D7").Value = "Left join case"
$ws.Rows.Item(7).RowHeight = 34

# --- Match the saved cursor/selection position recorded in the workbook. ---
[void]$ws.Range("C12").Select()
